{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" footer line, and the blank paragraph that\n// separates them from the \"LOM3206: Eletr\u00f4nica (Requisito)\" requirement\n// line above them (site footer/boilerplate pruned from the generated page).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the footer paragraphs by their exact text content so the script\n// is resilient to the surrounding content rather than relying on fixed\n// paragraph indices.\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (jupiterIndex === -1 && text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    jupiterIndex = i;\n  } else if (text.indexOf(\"Contact: luizeleno@usp.br\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (jupiterIndex !== -1 && copyrightIndex !== -1) {\n  // The empty paragraph right before \"Ver no Jupiter ...\" is also part of\n  // the removed block (as long as it immediately precedes it).\n  const blankIndex = jupiterIndex - 1;\n  if (blankIndex >= 0 && items[blankIndex].text === \"\") {\n    items[blankIndex].delete();\n  }\n  items[jupiterIndex].delete();\n  items[copyrightIndex].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"(c) 2020 . Contact: ...\" footer line, and the blank paragraph that\n# separates them from the \"LOM3206: Eletronica (Requisito)\" requirement\n# line above them (site footer/boilerplate pruned from the generated page).\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Locate the footer paragraphs by their exact text content so the script\n# is resilient to the surrounding content rather than relying on fixed\n# paragraph indices.\n$jupiterIndex = -1\n$copyrightIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($jupiterIndex -eq -1 -and $t -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        $jupiterIndex = $i\n    }\n    if ($t -like \"*Contact: luizeleno@usp.br*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($jupiterIndex -ne -1 -and $copyrightIndex -ne -1) {\n    # The empty paragraph right before \"Ver no Jupiter ...\" is also part of\n    # the removed block (as long as it immediately precedes it).\n    $blankIndex = $jupiterIndex - 1\n    $blankText = \"\"\n    if ($blankIndex -ge 1) {\n        $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd([char]13, [char]7)\n    }\n\n    # Delete starting from the highest index so earlier indices stay valid\n    # as paragraphs are removed.\n    $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n    if ($blankIndex -ge 1 -and $blankText -eq \"\") {\n        $d.Paragraphs.Item($blankIndex).Range.Delete()\n    }\n}\n"}
